$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.011.85"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +3.25%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.194.67"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +1.81%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "538.22"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.49%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.94"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +4.08%  "

$ws.Range("E8").Value = "  +2.60%  "

$ws.Range("E9").Value = "  +0.34%  "

$ws.Range("E10").Value = "  +4.15%  "

$ws.Range("E11").Value = "  +1.10%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.745.74"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +1.81%  "

$ws.Range("E13").Value = "  -1.45%  "

$ws.Range("E14").Value = "  +3.37%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.04"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.99%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "60.040.10"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +3.12%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.199.82"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +2.00%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.24"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.60%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.10"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.82%  "

$ws.Range("E20").Value = "  +1.71%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "383.01"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.54%  "

$ws.Range("E22").Value = "  +0.04%  "

$ws.Range("E23").Value = "  +2.62%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.32"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.20%  "

$ws.Range("E25").Value = "  +2.38%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.85"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +10.71%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.997"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.06%  "

$ws.Range("E28").Value = "  +2.41%  "

$ws.Range("E29").Value = "  +0.86%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.19"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.13%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.40"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +4.29%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.23"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +4.83%  "

$ws.Range("E34").Value = "  +5.74%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "156.47"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -3.27%  "

$ws.Range("E36").Value = "  -0.15%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.786.29"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +5.68%  "

$ws.Range("E38").Value = "  +0.61%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0713"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +5.41%  "

$ws.Range("E40").Value = "  +1.13%  "

$ws.Range("E41").Value = "  -0.11%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "39.79"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +2.07%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.730"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +4.19%  "

$ws.Range("E44").Value = "  +4.76%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.235.30"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.73%  "

$ws.Range("E46").Value = "  +2.73%  "

$ws.Range("E47").Value = "  +0.55%  "

$ws.Range("E48").Value = "  -0.90%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.805"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +7.43%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "20.56"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +1.18%  "

$ws.Range("E51").Value = "  +0.00%  "
